$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# Row 18: Estimation (hours) 2 -> 1
$ws.Cells.Item(18, 6).Value2 = 1

# Row 19: Priority Medium -> High
$ws.Cells.Item(19, 2).Value2 = "High"

# Row 21: Priority High -> Medium; Estimation (hours) blank -> 4
$ws.Cells.Item(21, 2).Value2 = "Medium"
$ws.Cells.Item(21, 6).Value2 = 4

# Row 22: new backlog entry
$ws.Cells.Item(22, 2).Value2 = "High"
$ws.Cells.Item(22, 3).Value2 = "To do"
$ws.Cells.Item(22, 4).Value2 = "to check or not the copyboard"
$ws.Cells.Item(22, 5).Value2 = "to avoid removing the search"
$ws.Cells.Item(22, 6).Value2 = 1

# Update the selected cell shown when the sheet is reopened
$ws.Activate()
$ws.Range("F20").Select()
